$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (so old B -> C)
$ws.Columns("B").Insert()
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# Update header row
$ws.Range("A1").Value = "ratingName_EN"
$ws.Range("B1").Value = "ratingName_CN"

# Copy style from A to B for header and data rows
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Range("A4").Copy()
$ws.Range("B4").PasteSpecial(-4122)

# Fill in Chinese values
$ws.Range("B2").Value = "Instructions_CN/ratingCS+1.png"
$ws.Range("B3").Value = "Instructions_CN/ratingCS+3.png"
$ws.Range("B4").Value = "Instructions_CN/ratingCS+4.png"

$ws.Range("C10").Select() | Out-Null
